# Add a new row for release/8.0.17 to the meta-sheet, following the same
# pattern as the existing release rows (branch name in col A, "X" placeholders
# for sit/uat/pre-prod/prod in cols B-E).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A20").Value = "release/8.0.17"
$ws.Range("B20").Value = "X"
$ws.Range("C20").Value = "X"
$ws.Range("D20").Value = "X"
$ws.Range("E20").Value = "X"

# Writing into row 20 picks up the column-level style (style="2" on <col>),
# but the previous last row (19) - and the target state - carry no explicit
# cell style. Reset to "Normal" so the new row matches that unstyled look.
$ws.Range("A20:E20").Style = "Normal"
